$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 83, shifting existing rows 83-111 down to 84-112
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new data entry
$ws.Range("A83").Value = 3
$ws.Range("B83").Value = "Femacal de La Calera"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 44463
$ws.Range("E83").Value = 5
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100101
$ws.Range("H83").Value = "Berries"
$ws.Range("I83").Value = 100112025
$ws.Range("J83").Value = "Frutilla"
$ws.Range("K83").Value = "Sin especificar"
$ws.Range("L83").Value = "Especial"
$ws.Range("M83").Value = 40
$ws.Range("N83").Value = 18000
$ws.Range("O83").Value = 18000
$ws.Range("P83").Value = 18000
$ws.Range("Q83").Value = "$/bandeja 7 kilos"
$ws.Range("R83").Value = "Provincia de Melipilla"
$ws.Range("S83").Value = 2571
$ws.Range("T83").Value = 7
